# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Overview: status for the a61deed4 row (row 3) failed handback
$ws1.Range("E3").Value = "Handback transform failed"
$ws1.Range("F3").Value = "Handback transform failed"

# zh-cn sheet: Status column (C3) for the a61deed4 row, widen Error Detail column,
# and set error message on row 3 (P3)
$ws2.Range("C3").Value = "Handback transform failed"
$ws2.Columns.Item(16).ColumnWidth = 39.1667
$ws2.Range("P3").Value = "Handback file name: gyjzona2.idt is different with handoff file name: a61deed4-0316-4cb0-9064-78c186258124.454ac6edf4b398f12466bde53f62fe08c8edbf94.zh-cn."

# de-de sheet: Status column (C3) for the a61deed4 row, widen Error Detail column,
# and set error message on row 3 (P3)
# (G2/J2/K2/G3 keep their existing handoff/handback values - unchanged by this edit)
$ws3.Range("C3").Value = "Handback transform failed"
$ws3.Columns.Item(16).ColumnWidth = 39.1667
$ws3.Range("P3").Value = "Handback file name: gyjzona2.idt is different with handoff file name: a61deed4-0316-4cb0-9064-78c186258124.454ac6edf4b398f12466bde53f62fe08c8edbf94.de-de."
